$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15; existing rows 15-80 shift down to 16-81
$ws.Rows.Item(15).Insert()

# Populate new row 15 with data (same as old row 15 except for the updated fields)
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 45030
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 100112010
$ws.Cells.Item(15, 7).Value = "Achicoria"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 10000
$ws.Cells.Item(15, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 556
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(15, 18).Value = "Hortaliza"
